$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update Target cluster (D2) from "FAPs" to "ECs", and stats ---
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.1847786666666667
$ws.Range("H2").Value = 0.5543360000000001
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.011963
$ws.Range("N2").Value = 0.035889
$ws.Range("O2").Value = 0.001886858915380773
$ws.Range("P2").Value = 0.001886858915380773
$ws.Range("Q2").Value = 0.002210507189333333
$ws.Range("R2").Value = 0.019894564704
$ws.Range("S2").Value = 0.001886858915380773
$ws.Range("T2").Value = 0.001886858915380773

# --- Row 3: update Target cluster (D3) from "sCs" to "FAPs", and stats ---
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.1847786666666667
$ws.Range("H3").Value = 0.5543360000000001
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.700797333333334
$ws.Range("N3").Value = 11.102392
$ws.Range("O3").Value = 0.5837066323177625
$ws.Range("P3").Value = 0.5837066323177625
$ws.Range("Q3").Value = 0.683828396856889
$ws.Range("R3").Value = 6.154455571712002
$ws.Range("S3").Value = 0.5837066323177625
$ws.Range("T3").Value = 0.5837066323177625

# --- Row 4: new row, Target cluster "sCs" ---
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Vip"
$ws.Range("C4").Value = "Adcyap1r1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.1847786666666667
$ws.Range("H4").Value = 0.5543360000000001
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.627406333333333
$ws.Range("N4").Value = 7.882218999999999
$ws.Range("O4").Value = 0.4144065087668568
$ws.Range("P4").Value = 0.4144065087668568
$ws.Range("Q4").Value = 0.4854886390648889
$ws.Range("R4").Value = 4.369397751584
$ws.Range("S4").Value = 0.4144065087668568
$ws.Range("T4").Value = 0.4144065087668568
